$d = $word.ActiveDocument

# 1. Remove the hidden _GoBack bookmark from its current location
#    ("Length of time to drive:" paragraph). It is re-added below, on
#    its own paragraph, further down after the new content.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Locate the "Ambient Temperature" paragraph - the new content is
#    inserted after the single blank paragraph that follows it.
$count = $d.Paragraphs.Count
$ambientIdx = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13,[char]7) -eq "Ambient Temperature") {
        $ambientIdx = $i
        break
    }
}

# 3. Of the two blank paragraphs following "Ambient Temperature", delete
#    the second one so only one remains (as in the target document).
$blank2 = $d.Paragraphs.Item($ambientIdx + 2)
$blank2.Range.Delete()

# 4. Insert the new paragraphs (TODO list items + "Simulation" section +
#    the relocated _GoBack bookmark) right after the remaining blank
#    paragraph. InsertXML always folds the very last <w:p> of the
#    fragment into the paragraph the insertion point sits in, so a
#    throwaway marker paragraph is appended after the bookmark
#    paragraph and then stripped out with Find/Replace, leaving the
#    bookmark on its own paragraph and "Elevation Profile..." untouched.
$blank1 = $d.Paragraphs.Item($ambientIdx + 1)
$insertPoint = $d.Range($blank1.Range.End, $blank1.Range.End)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Move buttons to panel</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Add elevation profile to </w:t></w:r><w:r><w:t>drive cycle</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Add gradient to constant speed and acceleration event</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Add length of time </w:t></w:r><w:r><w:t xml:space="preserve">panel </w:t></w:r><w:r><w:t>to run simulation (or run to % SOC remaining)</w:t></w:r></w:p><w:p/><w:p/><w:p><w:r><w:t>Add input warnings</w:t></w:r></w:p><w:p/><w:p/><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="ARRIVAL Apercu" w:eastAsia="Times New Roman" w:hAnsi="ARRIVAL Apercu" w:cs="Arial"/><w:b/><w:color w:val="222222"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:u w:val="single"/><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="ARRIVAL Apercu" w:eastAsia="Times New Roman" w:hAnsi="ARRIVAL Apercu" w:cs="Arial"/><w:b/><w:color w:val="222222"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:u w:val="single"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>Simulation</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="ARRIVAL Apercu" w:eastAsia="Times New Roman" w:hAnsi="ARRIVAL Apercu" w:cs="Arial"/><w:color w:val="222222"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="ARRIVAL Apercu" w:eastAsia="Times New Roman" w:hAnsi="ARRIVAL Apercu" w:cs="Arial"/><w:color w:val="222222"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>Side menu indicator to which section on</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="ARRIVAL Apercu" w:eastAsia="Times New Roman" w:hAnsi="ARRIVAL Apercu" w:cs="Arial"/><w:color w:val="222222"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="ARRIVAL Apercu" w:eastAsia="Times New Roman" w:hAnsi="ARRIVAL Apercu" w:cs="Arial"/><w:color w:val="222222"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>Make simulation results unfold when click simulate</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="ARRIVAL Apercu" w:eastAsia="Times New Roman" w:hAnsi="ARRIVAL Apercu" w:cs="Arial"/><w:color w:val="222222"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="ARRIVAL Apercu" w:eastAsia="Times New Roman" w:hAnsi="ARRIVAL Apercu" w:cs="Arial"/><w:color w:val="222222"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>Add other drive cycles and constant speed</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="ARRIVAL Apercu" w:eastAsia="Times New Roman" w:hAnsi="ARRIVAL Apercu" w:cs="Arial"/><w:color w:val="222222"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="ARRIVAL Apercu" w:eastAsia="Times New Roman" w:hAnsi="ARRIVAL Apercu" w:cs="Arial"/><w:color w:val="222222"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>Work out results screens</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:r><w:t>ZZPLACEHOLDERZZ</w:t></w:r></w:p>'
$insertPoint.InsertXML($xml)

$find = $d.Content.Find
$find.Execute("ZZPLACEHOLDERZZ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
